$wb = $excel.ActiveWorkbook

# --- Add handling for null cells -------------------------------------------
# The "DifferentTypes" sheet gains a new "Null" column (E) next to the
# existing Boolean/Number/String/Formula columns.
$ws = $wb.Worksheets.Item("DifferentTypes")
$ws.Range("E1").Value = "Null"

# This sheet becomes the active one, with the cursor parked just past the
# new data (as it would be right after typing the new header).
$ws.Activate() | Out-Null
$ws.Range("E6").Select() | Out-Null

# --- Housekeeping that accompanies the resave -------------------------------
# The workbook's external-connection defined names pick up a numeric suffix
# (as Excel does when it re-establishes the worksheet connections on save).
foreach ($n in $wb.Names) {
    $n.Name = $n.Name + "1"
}
